# Weekly update: insert a new daily price record for Pomelo (Fruta / Vega Central
# Mapocho de Santiago) at the top of the data table (row 5), pushing the
# existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5 (this shifts rows 5:49 down to 6:50 and
# extends the sheet dimension to A1:T50 automatically, inheriting the date
# number format from the row below for column D).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Cells.Item(5, 1).Value  = 9
$ws.Cells.Item(5, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value  = "Metropolitana"
$ws.Cells.Item(5, 4).Value  = 44462
$ws.Cells.Item(5, 5).Value  = 13
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100102
$ws.Cells.Item(5, 8).Value  = "Cítricos"
$ws.Cells.Item(5, 9).Value  = 100102006
$ws.Cells.Item(5, 10).Value = "Pomelo"
$ws.Cells.Item(5, 11).Value = "Start Ruby"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 550
$ws.Cells.Item(5, 14).Value = 7000
$ws.Cells.Item(5, 15).Value = 7500
$ws.Cells.Item(5, 16).Value = 7273
$ws.Cells.Item(5, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 520
$ws.Cells.Item(5, 20).Value = 14
